# Swap the presentation's design/theme colours: the deck originally used
# the "Integral" theme's colour set; the commit changes it over to the
# stock "Office Theme" colour set (fonts + effects/format scheme are
# already identical between the two themes, so only the 12-slot colour
# scheme on the slide master actually needs to change).
#
# Equivalent, in the UI, to opening Design > Variants > Colors and
# picking the different built-in colour set for the one-and-only slide
# master in this deck.

$p  = $ppt.ActivePresentation
$m  = $p.Slides.Item(1).Master
$cs = $m.ColorScheme

# Target RGB values straight from the "Office Theme" colour scheme
# (stored little-endian as R + G*256 + B*65536, matching the values the
# ColorScheme.Colors(n).RGB property already returns/accepts):
#   1  dk1      000000
#   2  lt1      FFFFFF
#   3  dk2      44546A
#   4  lt2      E7E6E6
#   5  accent1  5B9BD5
#   6  accent2  ED7D31
#   7  accent3  A5A5A5
#   8  accent4  FFC000
#   9  accent5  4472C4
#  10  accent6  70AD47
#  11  hlink    0563C1
#  12  folHlink 954F72
$cs.Colors(1).RGB  = 0
$cs.Colors(2).RGB  = 16777215
$cs.Colors(3).RGB  = 6968388
$cs.Colors(4).RGB  = 15132391
$cs.Colors(5).RGB  = 13998939
$cs.Colors(6).RGB  = 3243501
$cs.Colors(7).RGB  = 10855845
$cs.Colors(8).RGB  = 49407
$cs.Colors(9).RGB  = 12874308
$cs.Colors(10).RGB = 4697456
$cs.Colors(11).RGB = 12673797
$cs.Colors(12).RGB = 7491477
